# UKRI research council names changed
# - Remove the "IN-MCA-...AGRITECH-CATALYST..." / IPE Global row
# - Remove the "US-EIN-042103594-GCCI-3978870" / J-PAL row
# - Add a new row: GB-EDU-133903-PENDA / LSHTM / GB-GOV-1-300397
# - Keep autoFilter / _FilterDatabase ranges and sheet dimension in sync

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IATI activity IDs")

# --- Remove the two obsolete rows (found by their unique iati_id) ---
$target1 = $ws.Range("A:A").Find("IN-MCA-U74140DL1998PLC097579-1662_AGRITECH-CATALYST_DFID_UK")
if ($target1 -ne $null) {
    $ws.Rows($target1.Row).Delete()
}

$target2 = $ws.Range("A:A").Find("US-EIN-042103594-GCCI-3978870")
if ($target2 -ne $null) {
    $ws.Rows($target2.Row).Delete()
}

# --- Append the new row after the current last used row ---
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "GB-EDU-133903-PENDA"
$ws.Cells.Item($newRow, 2).Value = "Foreign, Commonwealth and Development Office"
$ws.Cells.Item($newRow, 3).Value = "GB-GOV-1-300397"
$ws.Cells.Item($newRow, 4).Value = "FCDO Research - Programmes"
$ws.Cells.Item($newRow, 5).Value = "LSHTM"

# --- Keep the AutoFilter range / sort-state references aligned with the shrunk data ---
$ws.AutoFilterMode = $false
$ws.Range("A1:E884").AutoFilter()

# --- Keep the hidden _FilterDatabase defined name aligned too ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "IATI activity IDs!_FilterDatabase") {
        $n.RefersTo = "='IATI activity IDs'!`$A`$1:`$E`$884"
    }
}
